$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows for the Guca2b -> Gucy2c LR-pair table (rows 2-10)
$rows = @(
  ,@("ECs", "Guca2b", "Gucy2c", "Inflammatory-Mac", 1.0, 0.5, 0.103181, 0.206362, 0.3503744070920869, 0.2893636473263989, 3.0, 1.0, 0.3282236666666666, 0.984671, 0.1589343097904449, 0.1589343097904449, 0.03386644615033333, 0.203198676902, 0.05568651455941721, 0.04598981156626693)
  ,@("ECs", "Guca2b", "Gucy2c", "Neutrophils", 1.0, 0.5, 0.103181, 0.206362, 0.3503744070920869, 0.2893636473263989, 3.0, 1.0, 1.716416, 5.149248, 0.8311326085767012, 0.8311326085767011, 0.177101519296, 1.062609115776, 0.2912075949449613, 0.2404995630296585)
  ,@("ECs", "Guca2b", "Gucy2c", "Resolving-Mac", 1.0, 0.5, 0.103181, 0.206362, 0.3503744070920869, 0.2893636473263989, 1.0, 0.3333333333333333, 0.02051333333333333, 0.06154, 0.009933081632854, 0.009933081632853998, 0.002116586246666666, 0.01269951748, 0.003480297587708519, 0.002874272730473495)
  ,@("FAPs", "Guca2b", "Gucy2c", "Inflammatory-Mac", 1.0, 0.3333333333333333, 0.1241823333333333, 0.372547, 0.4216891812734765, 0.5223905502006568, 3.0, 1.0, 0.3282236666666666, 0.984671, 0.1589343097904449, 0.1589343097904449, 0.04075958078188888, 0.366836227037, 0.0670208789717978, 0.08302578153719214)
  ,@("FAPs", "Guca2b", "Gucy2c", "Neutrophils", 1.0, 0.3333333333333333, 0.1241823333333333, 0.372547, 0.4216891812734765, 0.5223905502006568, 3.0, 1.0, 1.716416, 5.149248, 0.8311326085767012, 0.8311326085767011, 0.2131485438506667, 1.918336894656, 0.3504796292403979, 0.4341758206840901)
  ,@("FAPs", "Guca2b", "Gucy2c", "Resolving-Mac", 1.0, 0.3333333333333333, 0.1241823333333333, 0.372547, 0.4216891812734765, 0.5223905502006568, 1.0, 0.3333333333333333, 0.02051333333333333, 0.06154, 0.009933081632854, 0.009933081632853998, 0.002547393597777778, 0.02292654238, 0.00418867306128081, 0.005188947979374639)
  ,@("MuSCs", "Guca2b", "Gucy2c", "Inflammatory-Mac", 1.0, 0.5, 0.0671245, 0.134249, 0.2279364116344365, 0.1882458024729443, 3.0, 1.0, 0.3282236666666666, 0.984671, 0.1589343097904449, 0.1589343097904449, 0.02203184951316666, 0.132191097079, 0.03622691625922991, 0.02991871668698583)
  ,@("MuSCs", "Guca2b", "Gucy2c", "Neutrophils", 1.0, 0.5, 0.0671245, 0.134249, 0.2279364116344365, 0.1882458024729443, 3.0, 1.0, 1.716416, 5.149248, 0.8311326085767012, 0.8311326085767011, 0.115213565792, 0.691281394752, 0.189445384391342, 0.1564572248629526)
  ,@("MuSCs", "Guca2b", "Gucy2c", "Resolving-Mac", 1.0, 0.5, 0.0671245, 0.134249, 0.2279364116344365, 0.1882458024729443, 1.0, 0.3333333333333333, 0.02051333333333333, 0.06154, 0.009933081632854, 0.009933081632853998, 0.001376947243333333, 0.00826168346, 0.00226411098386467, 0.001869860923005865)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $row = $rows[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $row[$j]
  }
}
